# Refresh the crypto price table (Sheet1) with the latest scraped values.
# Source: GitHub Actions scheduled scrape, commit "Updated cryptos list on
# Wed Aug  9 23:42:22 UTC 2023 with GitHub Actions".
#
# Price (col D) and Volume/1h (col E) are stored as *text*, not numbers, in
# this sheet (t="inlineStr" in the OOXML) -- values like "1.854.94" are not
# valid numbers anyway, and even the ones that look numeric ("0.9989",
# "243.95", ...) must stay text to match the source data. Excel.Range.Value
# auto-converts a plain numeric-looking string to a Number, so those are
# written with a leading apostrophe (the normal Excel "force text" prefix),
# which Excel strips from the stored value but keeps the cell as Text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.585.91'
$ws.Range("E2").Value = '  -0.65%  '
# Row 3
$ws.Range("D3").Value = '1.855.15'
$ws.Range("E3").Value = '  -0.10%  '
# Row 4
$ws.Range("D4").Value = '''0.9989'
# Row 5
$ws.Range("D5").Value = '''243.95'
$ws.Range("E5").Value = '  -0.50%  '
# Row 6
$ws.Range("D6").Value = '''0.6451'
$ws.Range("E6").Value = '  +0.68%  '
# Row 7
$ws.Range("D7").Value = '''0.9995'
$ws.Range("E7").Value = '  -0.07%  '
# Row 8
$ws.Range("E8").Value = '  +1.44%  '
# Row 9
$ws.Range("D9").Value = '''0.07536'
$ws.Range("E9").Value = '  +0.49%  '
# Row 10
$ws.Range("D10").Value = '''24.41'
$ws.Range("E10").Value = '  +0.93%  '
# Row 11
$ws.Range("D11").Value = '''0.07671'
$ws.Range("E11").Value = '  -0.17%  '
# Row 12
$ws.Range("D12").Value = '1.911.13'
$ws.Range("E12").Value = '  +2.60%  '
# Row 13
$ws.Range("D13").Value = '''5.054'
$ws.Range("E13").Value = '  -0.11%  '
# Row 14
$ws.Range("D14").Value = '''0.6908'
$ws.Range("E14").Value = '  +0.82%  '
# Row 15
$ws.Range("D15").Value = '''84.02'
$ws.Range("E15").Value = '  +0.04%  '
# Row 16
$ws.Range("D16").Value = '''0.000009597'
$ws.Range("E16").Value = '  +1.02%  '
# Row 17
$ws.Range("D17").Value = '''6.268'
$ws.Range("E17").Value = '  +3.02%  '
# Row 18
$ws.Range("D18").Value = '2.164.79'
$ws.Range("E18").Value = '  +2.37%  '
# Row 19
$ws.Range("D19").Value = '29.607.60'
$ws.Range("E19").Value = '  -0.53%  '
# Row 20
$ws.Range("D20").Value = '''237.44'
$ws.Range("E20").Value = '  -0.94%  '
# Row 21
$ws.Range("D21").Value = '''12.63'
$ws.Range("E21").Value = '  -0.40%  '
# Row 22
$ws.Range("D22").Value = '''0.9999'
$ws.Range("E22").Value = '  +0.01%  '
# Row 23
$ws.Range("D23").Value = '''7.725'
$ws.Range("E23").Value = '  +4.01%  '
# Row 24
$ws.Range("D24").Value = '''0.9999'
$ws.Range("E24").Value = '  -0.12%  '
# Row 25
$ws.Range("D25").Value = '''157.33'
$ws.Range("E25").Value = '  -0.99%  '
# Row 26
$ws.Range("D26").Value = '''0.1418'
$ws.Range("E26").Value = '  -1.05%  '
# Row 27
$ws.Range("D27").Value = '''8.533'
$ws.Range("E27").Value = '  +0.13%  '
# Row 28
$ws.Range("D28").Value = '''17.84'
$ws.Range("E28").Value = '  -0.61%  '
# Row 29
$ws.Range("D29").Value = '''1.493'
$ws.Range("E29").Value = '  -0.66%  '
# Row 30
$ws.Range("E30").Value = '  -3.37%  '
# Row 31
$ws.Range("E31").Value = '  -1.37%  '
# Row 32
$ws.Range("D32").Value = '''4.143'
$ws.Range("E32").Value = '  -0.09%  '
# Row 33
$ws.Range("D33").Value = '''4.085'
$ws.Range("E33").Value = '  -0.54%  '
# Row 34
$ws.Range("D34").Value = '''1.886'
$ws.Range("E34").Value = '  +0.49%  '
# Row 35
$ws.Range("D35").Value = '''1.177'
$ws.Range("E35").Value = '  +1.44%  '
# Row 36
$ws.Range("D36").Value = '''0.7241'
$ws.Range("E36").Value = '  -1.04%  '
# Row 37
$ws.Range("D37").Value = '''2.605'
$ws.Range("E37").Value = '  +0.01%  '
# Row 38
$ws.Range("D38").Value = '''2.789'
$ws.Range("E38").Value = '  -2.38%  '
# Row 39
$ws.Range("D39").Value = '''0.01779'
$ws.Range("E39").Value = '  -0.93%  '
# Row 40
$ws.Range("D40").Value = '1.211.00'
$ws.Range("E40").Value = '  -0.14%  '
# Row 41
$ws.Range("D41").Value = '''0.9131'
$ws.Range("E41").Value = '  -1.22%  '
# Row 42
$ws.Range("D42").Value = '''6.189'
$ws.Range("E42").Value = '  +0.09%  '
# Row 43
$ws.Range("D43").Value = '2.072.06'
$ws.Range("E43").Value = '  +2.33%  '
# Row 44
$ws.Range("D44").Value = '''0.9995'
$ws.Range("E44").Value = '  -0.11%  '
# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '''0.00000000128'
$ws.Range("E45").Value = '  +5.97%  '
# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''102.05'
$ws.Range("E46").Value = '  +0.01%  '
# Row 47
$ws.Range("D47").Value = '''67.33'
$ws.Range("E47").Value = '  +1.46%  '
# Row 48
$ws.Range("D48").Value = '''7.372'
$ws.Range("E48").Value = '  +10.07%  '
# Row 49
$ws.Range("D49").Value = '''0.4071'
$ws.Range("E49").Value = '  -0.27%  '
# Row 50
$ws.Range("D50").Value = '''9.193'
$ws.Range("E50").Value = '  -1.67%  '
# Row 51
$ws.Range("D51").Value = '''1.668'
$ws.Range("E51").Value = '  +2.33%  '
